$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Relabel the header row: "<Name>_old" -> "<Name>_FV2410", "<Name>_new" -> "<Name>_FV2504"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -like "*_old") {
        $cell.Value = $val -replace "_old$", "_FV2410"
    } elseif ($val -like "*_new") {
        $cell.Value = $val -replace "_new$", "_FV2504"
    }
}

# 2. Turn the data range A1:U62 into a proper Excel Table ("Table1")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U62"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# 3. Freeze the header row (pane split below row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
